# "Added New Mac-Address and Document Types"
# Appends one new test-data row (regcntr_id/usr_id/lang_code/is_active/cr_by/cr_dtimes)
# to the bottom of the master-reg_center_user sheet, matching the existing rows'
# pattern, and updates the sheet's view/selection to the new bottom of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count + 1, 1).Row

$ws.Cells.Item($newRow, 1).Value = 10002        # regcntr_id
$ws.Cells.Item($newRow, 2).Value = 110032       # usr_id
$ws.Cells.Item($newRow, 3).Value = "eng"        # lang_code
$ws.Cells.Item($newRow, 4).Value = $true        # is_active
$ws.Cells.Item($newRow, 5).Value = "superadmin" # cr_by
$ws.Cells.Item($newRow, 6).Value = "now()"      # cr_dtimes

# Scroll the view down and move the selection near the newly added row,
# mirroring the author's updated sheetView/selection.
$ws.Range("C31").Select()
$excel.ActiveWindow.ScrollRow = 22
